# Update the "tijdsbesteding" workbook: add a new week row (row 29) with its
# hours entry, refresh the totals/averages that now span the new row, and
# reset the sheet view (scroll position + selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new week row --------------------------------------------------
$ws.Range("A29").Value = "week 16-22/05/2016"
$ws.Range("B29").Formula = "=6+2+2"

# Copy the formatting from the row above (row 28) onto the new row so the
# new cells pick up the same styles (label fill + currency number format).
$ws.Range("A28:B28").Copy()
$ws.Range("A29:B29").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Recalculate so SUM/AVERAGE/ratio formulas pick up the new row 29 value.
$excel.CalculateFull()

# --- Reset the view: scroll back to the top and select E2 ------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("E2").Select()
